# Update "想去人数" (want-to-go count) figures across sheets to match the
# latest scrape output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 9370
$ws1.Range("F9").Value  = 5725
$ws1.Range("F12").Value = 6597
$ws1.Range("F14").Value = 461
$ws1.Range("F24").Value = 10664
$ws1.Range("F28").Value = 2501
$ws1.Range("F37").Value = 1471
$ws1.Range("F39").Value = 5449
$ws1.Range("F47").Value = 1427
$ws1.Range("F49").Value = 1109

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F20").Value = 20

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 76

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 9370
$ws4.Range("F12").Value = 5725
$ws4.Range("F13").Value = 5725
$ws4.Range("F16").Value = 6597
$ws4.Range("F17").Value = 6597
$ws4.Range("F19").Value = 461
$ws4.Range("F28").Value = 10664
$ws4.Range("F32").Value = 2501
$ws4.Range("F38").Value = 1471
$ws4.Range("F40").Value = 5449
$ws4.Range("F41").Value = 20
$ws4.Range("F49").Value = 1427
$ws4.Range("F51").Value = 1109
